$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing value for F11 (continuing the data series)
$ws.Range("F11").Value = 340

# Update the selection to reflect where editing left off
$ws.Range("F11").Select()
